$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @(0.6117092866756392, 1690369242.192249, 0.9700400031215803, 1690369242.259881, 0.06763243675231934)
    3  = @(0.2893674293405114, 1690369247.180922, 0.7457110641474254, 1690369247.249108, 0.06818628311157227)
    4  = @(0.4703903095558546, 1690369252.184102, 0.5088808143078496, 1690369252.240347, 0.05624556541442871)
    5  = @(0.6117092866756392, 1690369257.198658, 0.2640038046884596, 1690369257.239956, 0.04129815101623535)
    6  = @(0.4966352624495289, 1690369262.184294, 0.9556950817962525, 1690369262.24972,  0.06542539596557617)
    7  = @(0.4175639300134589, 1690369267.17212,  0.7199734245730661, 1690369267.240569, 0.06844902038574219)
    8  = @(0.3024899057873485, 1690369272.17943,  0.4483742410832863, 1690369272.239869, 0.06043910980224609)
    9  = @(0.4801480484522207, 1690369277.203125, 0.1872779038191648, 1690369277.229915, 0.02678990364074707)
    10 = @(0.595222072678331,  1690369282.192439, 0.9001615732764729, 1690369282.239987, 0.04754829406738281)
    11 = @(0.5262449528936743, 1690369287.1872,   0.5847835091905227, 1690369287.239927, 0.05272698402404785)
    12 = @(0.4407806191117093, 1690369292.191486, 0.3247972226800773, 1690369292.229494, 0.03800797462463379)
    13 = @(0.335464333781965,  1690369297.183312, 0.9859821680308416, 1690369297.240759, 0.05744743347167969)
}

foreach ($row in $data.Keys) {
    $values = $data[$row]
    $ws.Cells.Item($row, 1).Value = $values[0]
    $ws.Cells.Item($row, 2).Value = $values[1]
    $ws.Cells.Item($row, 3).Value = $values[2]
    $ws.Cells.Item($row, 4).Value = $values[3]
    $ws.Cells.Item($row, 5).Value = $values[4]
}
